# Fix to manufacturing files: update Rotation values (column E) for
# several components and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Rotation (column E) values
$ws.Range("E7").Value = 0
$ws.Range("E13").Value = 180
$ws.Range("E14").Value = 180
$ws.Range("E46").Value = 90
$ws.Range("E48").Value = 90
$ws.Range("E49").Value = 90

# Update the active cell / selection to match the saved view state
$ws.Range("E14").Select()
